$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell (H1) to the new headers so
# they match the bold/bordered/centered look of the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data for the new I (I0) and J (IF) columns, rows 2-19
$iValues = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,7,7)
$jValues = @(5,5,6,5,3,4,4,5,5,3,5,6,5,5,5,5,8,7)

for ($r = 0; $r -lt 18; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
